$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("G2").Value = "2016-09-06 14:40:59"

$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("H2").Value = "2016-09-06 14:40:39"
$ws2.Range("K2").Value = "2016-09-06 14:42:21"

$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("H2").Value = "2016-09-06 14:40:59"
$ws3.Range("K2").Value = "2016-09-06 14:42:47"
